$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = -123.47338606118794;  C = 0.09136807702319806;  D = 62.406021743 }
    3  = @{ B = -117.61674244158768;  C = 0.0646728753907774;   D = 78.841764898 }
    4  = @{ B = -122.00142694489112;  C = 0.07604737257360063;  D = 67.70893263 }
    5  = @{ B = -121.30725095787741;  C = 0.09996849057912861;  D = 95.756308402 }
    6  = @{ B = -120.12292887523154;  C = 0.0840953498554457;   D = 69.010277998 }
    7  = @{ B = -118.89272648754448;  C = 0.08608142800931605;  D = 62.188166612 }
    8  = @{ B = -115.03704432900317;  C = 0.0;                  D = 77.680880971 }
    9  = @{ B = -119.44147912757283;  C = 0.09684089792055087;  D = 96.448602382 }
    10 = @{ B = -120.14729823432504;  C = 0.09943172612207225;  D = 72.967578873 }
    11 = @{ B = -116.2650930538125;   C = 0.0936922930803802;   D = 99.783351422 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
}
